$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The published dataset gained one additional weekly record. It is inserted
# as the new row 109 (pushing the former rows 109-176 down to 110-177), so
# the sheet's used range grows from A1:R176 to A1:R177.
$ws.Rows.Item(109).Insert()

$ws.Range("A109").Value = 5
$ws.Range("B109").Value = "Macroferia Regional de Talca"
$ws.Range("C109").Value = "Maule"
$ws.Range("D109").Value = 44438
$ws.Range("E109").Value = 7
$ws.Range("F109").Value = 100114013
$ws.Range("G109").Value = "Zanahoria"
$ws.Range("H109").Value = "Sin especificar"
$ws.Range("I109").Value = "Primera"
$ws.Range("J109").Value = 500
$ws.Range("K109").Value = 6000
$ws.Range("L109").Value = 6000
$ws.Range("M109").Value = 6000
$ws.Range("N109").Value = "$/saco 20 kilos"
$ws.Range("O109").Value = "Región de La Araucanía"
$ws.Range("P109").Value = 300
$ws.Range("Q109").Value = 20
$ws.Range("R109").Value = "Hortaliza"

# Keep the date column's display format consistent with the rest of the
# column (matches the style already inherited from the row-insert, but set
# explicitly so it is not dependent on that behaviour).
$ws.Range("D109").NumberFormat = "YYYY-MM-DD HH:MM:SS"
